$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6087971329689026
$ws.Range("B1").Value = 1.360681056976318
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.449792623519897
$ws.Range("E1").Value = 1.395982503890991
